$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: create the new row 29 by duplicating row 28 (all 18 columns: A..R) ---
# This grabs the "static" columns (A,B,C,E,F,G,H,I,N,O,Q,R) as well as the
# "varying" columns (D,J,K,L,M,P) which will be overwritten with their final
# values in step 2 below. Also copy the Fecha (D) column number format so the
# new row matches the date-styled cells above it.
for ($col = 1; $col -le 18; $col++) {
    $ws.Cells.Item(29, $col).Value2 = $ws.Cells.Item(28, $col).Value2
}
$ws.Cells.Item(29, 4).NumberFormat = $ws.Cells.Item(28, 4).NumberFormat

# --- Step 2: apply the new Fecha/Volumen/Precio values for rows 2..29 ---
# (each row effectively "shifted down" one slot, with a brand new entry in row 2
# and the former row 28 entry now living in the new row 29)
$rows = @(
    @{Row=2;  D=44756; J=120; K=3300; L=3300; M=3300; P=3300},
    @{Row=3;  D=44452; J=120; K=2300; L=2300; M=2300; P=2300},
    @{Row=4;  D=44474; J=20;  K=1600; L=1600; M=1600; P=1600},
    @{Row=5;  D=44740; J=50;  K=2500; L=2500; M=2500; P=2500},
    @{Row=6;  D=44483; J=50;  K=2200; L=2200; M=2200; P=2200},
    @{Row=7;  D=44447; J=75;  K=2200; L=2200; M=2200; P=2200},
    @{Row=8;  D=44749; J=80;  K=2500; L=2500; M=2500; P=2500},
    @{Row=9;  D=44669; J=60;  K=6250; L=6250; M=6250; P=6250},
    @{Row=10; D=44484; J=40;  K=2200; L=2200; M=2200; P=2200},
    @{Row=11; D=44706; J=90;  K=4700; L=4700; M=4700; P=4700},
    @{Row=12; D=44677; J=20;  K=5500; L=5500; M=5500; P=5500},
    @{Row=13; D=44741; J=100; K=2500; L=2500; M=2500; P=2500},
    @{Row=14; D=44496; J=40;  K=2200; L=2200; M=2200; P=2200},
    @{Row=15; D=44679; J=30;  K=5500; L=5500; M=5500; P=5500},
    @{Row=16; D=44720; J=100; K=3600; L=3600; M=3600; P=3600},
    @{Row=17; D=44203; J=30;  K=2000; L=2000; M=2000; P=2000},
    @{Row=18; D=44497; J=50;  K=2200; L=2200; M=2200; P=2200},
    @{Row=19; D=44707; J=100; K=4700; L=4700; M=4700; P=4700},
    @{Row=20; D=44685; J=60;  K=5000; L=6000; M=5333; P=5333},
    @{Row=21; D=44487; J=50;  K=2200; L=2200; M=2200; P=2200},
    @{Row=22; D=44755; J=90;  K=3300; L=3300; M=3300; P=3300},
    @{Row=23; D=44476; J=30;  K=2200; L=2200; M=2200; P=2200},
    @{Row=24; D=44747; J=80;  K=2500; L=2500; M=2500; P=2500},
    @{Row=25; D=44453; J=20;  K=2300; L=2300; M=2300; P=2300},
    @{Row=26; D=44754; J=50;  K=3300; L=3300; M=3300; P=3300},
    @{Row=27; D=44719; J=80;  K=3600; L=3600; M=3600; P=3600},
    @{Row=28; D=44473; J=140; K=1600; L=1600; M=1600; P=1600},
    @{Row=29; D=44753; J=130; K=2700; L=3300; M=2931; P=2931}
)

foreach ($entry in $rows) {
    $r = $entry.Row
    $ws.Cells.Item($r, 4).Value2  = $entry.D   # Fecha
    $ws.Cells.Item($r, 10).Value2 = $entry.J   # Volumen
    $ws.Cells.Item($r, 11).Value2 = $entry.K   # Precio minimo
    $ws.Cells.Item($r, 12).Value2 = $entry.L   # Precio maximo
    $ws.Cells.Item($r, 13).Value2 = $entry.M   # Precio promedio ponderado
    $ws.Cells.Item($r, 16).Value2 = $entry.P   # Precio $/Kg
}
